# The review-dates list is auto-generated from the document repository.
# On regeneration, four PDFs that were removed / superseded in the source
# folder (old SARI guidance docs) dropped out of the list entirely.
# Delete their rows here; Excel will re-pack sharedStrings on save so the
# remaining "File"/"Review date" pairs keep their original order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filesToRemove = @(
    "Infection_and_sepsis/SARI/Reporting influenza deaths.pdf",
    "Infection_and_sepsis/SARI/Setup Guide for Jupiter Hoods.pdf",
    "Infection_and_sepsis/SARI/Management of patients with SARI-additional Information.pdf",
    "Infection_and_sepsis/SARI/Management of Patients with severe acute respiratory infection SARI.pdf"
)

$lastRow = $ws.UsedRange.Rows.Count()

# Walk bottom-to-top so deleting a row never shifts the row index of one
# not yet visited.
for ($i = $lastRow; $i -ge 2; $i--) {
    $name = $ws.Cells.Item($i, 1).Value()
    if ($filesToRemove -contains $name) {
        $ws.Rows.Item($i).Delete()
    }
}
